$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: WE Fashion row - only the source URL (ProductsOverview URL) changes ---
$ws.Range("I2").Value = "https://www.wefashion.nl/nl_NL/outlet/men"

# --- Row 3: replace "Nasty gal"/isawitfirst.com selectors with "Free people" selectors ---
$ws.Range("A3").Value = "Free people"
$ws.Range("B3").Value = ".c-pwa-tile-view-outer"
$ws.Range("C3").Value = ".c-pwa-product-tile"
$ws.Range("D3").Value = ".c-pwa-product-meta-heading"
$ws.Range("E3").Value = ".c-pwa-product-price__current"
$ws.Range("F3").Value = ".c-pwa-product-price__original"
$ws.Range("G3").Value = ".c-pwa-image-viewer__img:src"
$ws.Range("H3").Value = ".c-pwa-product-tile__link"
$ws.Range("I3").Value = "https://www.freepeople.com/sale-all/?currency=EUR"
# Plain `.Value = "FALSE"` gets auto-typed to a real Excel boolean; the
# source sheet stores this column as literal text, so build it as a text
# formula result and flatten it to a static value (keeps it t="s" without
# leaving a quotePrefix artifact behind).
$ws.Range("J3").Formula = "=""FALSE"""
$ws.Range("J3").Copy()
$ws.Range("J3").PasteSpecial(-4163)

# --- Row 4: brand new "Forever21" row ---
$ws.Range("A4").Value = "Forever21"
$ws.Range("B4").Value = ".product-grid"
$ws.Range("C4").Value = ".product-tile"
$ws.Range("D4").Value = ".pdp__name"
$ws.Range("E4").Value = ".price__default--discount"
$ws.Range("F4").Value = ".price__original > .value"
$ws.Range("G4").Value = ".product-gallery__img:src"
$ws.Range("H4").Value = ".product-tile__anchor"
$ws.Range("I4").Value = "https://www.forever21.com/us/shop/catalog/category/21men/mens-sale"
$ws.Range("J4").Formula = "=""TRUE"""
$ws.Range("J4").Copy()
$ws.Range("J4").PasteSpecial(-4163)

# --- Hyperlinks: the host only lets us delete ALL hyperlinks at once, so clear
# them and recreate the full, final set (I2, I3, I4) in order so the
# generated relationship ids line up as rId1/rId2/rId3. ---
$ws.Range("I2").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("I2"), "https://www.wefashion.nl/nl_NL/outlet/men")
$ws.Hyperlinks.Add($ws.Range("I3"), "https://www.freepeople.com/sale-all/?currency=EUR")
$ws.Hyperlinks.Add($ws.Range("I4"), "https://www.forever21.com/us/shop/catalog/category/21men/mens-sale")

$ws.Range("I2").Style = "Hyperlink"
$ws.Range("I3").Style = "Hyperlink"
$ws.Range("I4").Style = "Hyperlink"

# --- Selection moved from J4 to G5 (and the view no longer needs to be
# scrolled to keep column J in view) ---
$ws.Range("G5").Select()
